# Update the "target radius" coordinates on Sheet1 from 16 -> 12
# (diagonal values correspondingly from 11.314 -> 8.485), leaving the
# "on/off" flag in column C untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @(12, 0)
    3  = @(0, 12)
    4  = @(8.485, 8.485)
    5  = @(-8.485, 8.485)
    6  = @(-12, 0)
    7  = @(0, -12)
    8  = @(-8.485, -8.485)
    9  = @(8.485, -8.485)
    10 = @(8.485, 8.485)
    11 = @(12, 0)
    12 = @(0, 12)
    13 = @(-12, 0)
    14 = @(-8.485, 8.485)
    15 = @(-8.485, -8.485)
    16 = @(0, -12)
    17 = @(8.485, -8.485)
    18 = @(0, 12)
    19 = @(8.485, 8.485)
    20 = @(12, 0)
    21 = @(-8.485, 8.485)
    22 = @(-12, 0)
    23 = @(-8.485, -8.485)
    24 = @(8.485, -8.485)
    25 = @(0, -12)
    26 = @(-8.485, 8.485)
    27 = @(8.485, 8.485)
    28 = @(0, 12)
    29 = @(0, -12)
    30 = @(-12, 0)
    31 = @(-8.485, -8.485)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}

# Update the active selection from F8 to E8, matching the saved view state.
$ws.Range("E8").Select()
